$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-redundant "municipio / nº de casos / nº de óbitos" header row (row 2).
# This shifts every data row up by one (old row 3 "aruja" becomes new row 2, etc.).
$ws.Rows.Item(2).Delete()

# Remove the trailing "(vazio)" row, which after the shift above is now row 51.
$ws.Rows.Item(51).Delete()
